$wb = $excel.ActiveWorkbook

$wsTeacher = $wb.Worksheets.Item("Teacher")
$wsClasses = $wb.Worksheets.Item("Classes")

# --- Classes sheet data edits ---
# Row 2: class overlap note for MATH 313, and first Num_Sections value
$wsClasses.Range("C2").Value = "MATH 401, MATH 125"
$wsClasses.Range("F2").Value = 1

# Row 3: Num_Sections for MATH 401
$wsClasses.Range("F3").Value = 1

# Row 4: Num_Sections for MATH 125
$wsClasses.Range("F4").Value = 1

# Row 5: class overlap note for MATH 130, and Num_Sections
$wsClasses.Range("C5").Value = "STAT 102"
$wsClasses.Range("F5").Value = 1

# Row 6: Num_Sections for STAT 102
$wsClasses.Range("F6").Value = 2

# --- Selection / active sheet changes ---
# Teacher tab loses the selected/active state, selection moves to B3
$wsTeacher.Select() | Out-Null
$wsTeacher.Range("B3").Select() | Out-Null

# Classes tab becomes the active/selected tab, selection moves to F7
$wsClasses.Select() | Out-Null
$wsClasses.Range("F7").Select() | Out-Null
